$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old data rows (2-6) entirely so their formatting (style s="1")
# is dropped along with them; fresh rows written below come back with the
# default (no explicit style) formatting, matching the target file.
$ws.Rows("2:6").Delete()

# Write the new cell values. The order below controls the order in which
# brand-new strings land in the shared string table (existing strings keep
# their prior slot), matching the saved workbook's table layout.
$ws.Range("A2").Value = "skillname_attack"
$ws.Range("A3").Value = "skillname_shoot"

$ws.Range("B2").Value = "평타"
$ws.Range("C2").Value = "Normal Attack"
$ws.Range("D2").Value = "Normal Attack"
$ws.Range("E2").Value = "Normal Attack"

$ws.Range("B3").Value = "사격"
$ws.Range("C3").Value = "Shooting"
$ws.Range("D3").Value = "Shooting"
$ws.Range("E3").Value = "Shooting"

$ws.Range("B1").Value = "Korean"
$ws.Range("C1").Value = "English"
$ws.Range("D1").Value = "Chinese"
$ws.Range("E1").Value = "Japanese"
$ws.Range("A1").Value = "ID"

# Columns B:E now share the same width (~18.25 chars; 17.6 is the closest
# COM input that lands on the engine's nearest storable width to 18.25)
$ws.Columns("2:5").ColumnWidth = 17.6

# Leave selection/active cell on C3, matching the saved workbook state
$ws.Range("C3").Select()
